# appdev_project2/doc/db.xlsx — "Added Bid and fixed class relationships"
#
# The ER-diagram sheet modeled a many-to-one "JobOwners" table (jobId -> userId)
# separately from the Job table. This edit folds that relationship directly
# into the Job table as a new "owner" column, removing the now-redundant
# JobOwners table, and the Bid/Role tables below it shift up to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fold the JobOwners (jobId -> userId) mapping into the Job table as a new
#    "owner" column (G), using the same owners the JobOwners table listed:
#    J1->U2, J2->U2, J3->U1, J4->U3.
$ws.Cells.Item(9, 7).Value = "owner"
$ws.Cells.Item(10, 7).Value = "U2"
$ws.Cells.Item(11, 7).Value = "U2"
$ws.Cells.Item(12, 7).Value = "U1"
$ws.Cells.Item(13, 7).Value = "U3"

# Give the new header cell (G8) the same bottom-border treatment as the rest
# of the Job header row (A8:F8), so the header rule spans the new column too.
$border = $ws.Cells.Item(8, 7).Borders.Item(9)
$border.LineStyle = 1
$border.Weight = 2

# Column G now holds short userId codes (U1..U4) instead of the old wider
# "owner name" style content — shrink it to fit.
$ws.Columns.Item(7).ColumnWidth = 6.7109375

# 2) The JobOwners table (header + sub-header + 4 rows, plus its trailing
#    blank spacer row) is now redundant — delete rows 15:21 outright. This
#    shifts the Bid table (was rows 22:29) up to rows 15:22, and the Role
#    table (was rows 31:37) up to rows 24:30.
$ws.Range("A15:A21").EntireRow.Delete()

# 3) Restore the originally-selected cell now that the layout has shifted.
$ws.Range("B19").Select()
